$wb = $excel.ActiveWorkbook

# ===== Sheet 1: sheet1 =====
$ws = $wb.Worksheets.Item(1)

$ws.Range('A2').Value = 'Última actualización: 16:52:47'
$ws.Range('A3').Value = 'Total filas: 347'

# Update rows whose values changed due to re-sort
$ws.Cells.Item(47,1).Value = '05:49:10'
$ws.Cells.Item(47,2).Value = '07:32'
$ws.Cells.Item(47,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(47,4).Value = 103
$ws.Cells.Item(47,5).Value = 'LP1912'
$ws.Cells.Item(49,1).Value = '05:49:10'
$ws.Cells.Item(49,2).Value = '07:32'
$ws.Cells.Item(49,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(49,4).Value = 103
$ws.Cells.Item(49,5).Value = 'LP1912'
$ws.Cells.Item(80,1).Value = '07:19:37'
$ws.Cells.Item(80,2).Value = '08:43'
$ws.Cells.Item(80,3).Value = '14_ABASTO'
$ws.Cells.Item(80,4).Value = 84
$ws.Cells.Item(80,5).Value = 'LP1912'
$ws.Cells.Item(81,1).Value = '08:19:33'
$ws.Cells.Item(81,2).Value = '08:43'
$ws.Cells.Item(81,3).Value = '16_SANTA ANA'
$ws.Cells.Item(81,4).Value = 24
$ws.Cells.Item(81,5).Value = 'LP1912'
$ws.Cells.Item(89,1).Value = '08:19:33'
$ws.Cells.Item(89,2).Value = '09:02'
$ws.Cells.Item(89,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(89,4).Value = 43
$ws.Cells.Item(89,5).Value = 'LP1912'
$ws.Cells.Item(90,1).Value = '08:48:09'
$ws.Cells.Item(90,2).Value = '09:02'
$ws.Cells.Item(90,3).Value = '215A_EL PATO'
$ws.Cells.Item(90,4).Value = 14
$ws.Cells.Item(90,5).Value = 'LP1912'
$ws.Cells.Item(124,1).Value = '08:19:33'
$ws.Cells.Item(124,2).Value = '10:12'
$ws.Cells.Item(124,3).Value = '15_ABASTO'
$ws.Cells.Item(124,4).Value = 113
$ws.Cells.Item(124,5).Value = 'LP1912'
$ws.Cells.Item(125,1).Value = '09:25:56'
$ws.Cells.Item(125,2).Value = '10:12'
$ws.Cells.Item(125,3).Value = '10_OLMOS'
$ws.Cells.Item(125,4).Value = 47
$ws.Cells.Item(125,5).Value = 'LP1912'
$ws.Cells.Item(156,1).Value = '11:17:08'
$ws.Cells.Item(156,2).Value = '11:17'
$ws.Cells.Item(156,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(156,4).Value = 0
$ws.Cells.Item(156,5).Value = 'LP1912'
$ws.Cells.Item(157,1).Value = '11:17:08'
$ws.Cells.Item(157,2).Value = '11:17'
$ws.Cells.Item(157,3).Value = '16_SANTA ANA'
$ws.Cells.Item(157,4).Value = 0
$ws.Cells.Item(157,5).Value = 'LP1912'
$ws.Cells.Item(243,1).Value = '14:00:52'
$ws.Cells.Item(243,2).Value = '14:04'
$ws.Cells.Item(243,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(243,4).Value = 4
$ws.Cells.Item(243,5).Value = 'LP1912'
$ws.Cells.Item(244,1).Value = '12:27:08'
$ws.Cells.Item(244,2).Value = '14:04'
$ws.Cells.Item(244,3).Value = '17_ROMERO'
$ws.Cells.Item(244,4).Value = 97
$ws.Cells.Item(244,5).Value = 'LP1912'
$ws.Cells.Item(288,1).Value = '15:51:48'
$ws.Cells.Item(288,2).Value = '16:02'
$ws.Cells.Item(288,3).Value = '16_SANTA ANA'
$ws.Cells.Item(288,4).Value = 11
$ws.Cells.Item(288,5).Value = 'LP1912'
$ws.Cells.Item(289,1).Value = '14:44:25'
$ws.Cells.Item(289,2).Value = '16:02'
$ws.Cells.Item(289,3).Value = '27_EL RETIRO'
$ws.Cells.Item(289,4).Value = 78
$ws.Cells.Item(289,5).Value = 'LP1912'
$ws.Cells.Item(309,1).Value = '16:18:00'
$ws.Cells.Item(309,2).Value = '16:43'
$ws.Cells.Item(309,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(309,4).Value = 25
$ws.Cells.Item(309,5).Value = 'LP1912'
$ws.Cells.Item(310,1).Value = '14:44:25'
$ws.Cells.Item(310,2).Value = '16:43'
$ws.Cells.Item(310,3).Value = '225_GOMEZ'
$ws.Cells.Item(310,4).Value = 119
$ws.Cells.Item(310,5).Value = 'LP1912'
$ws.Cells.Item(314,1).Value = '16:52:47'
$ws.Cells.Item(314,2).Value = '16:53'
$ws.Cells.Item(314,3).Value = '16_SANTA ANA'
$ws.Cells.Item(314,4).Value = 1
$ws.Cells.Item(314,5).Value = 'LP1912'
$ws.Cells.Item(315,1).Value = '16:52:47'
$ws.Cells.Item(315,2).Value = '16:53'
$ws.Cells.Item(315,3).Value = '15_ABASTO'
$ws.Cells.Item(315,4).Value = 1
$ws.Cells.Item(315,5).Value = 'LP1912'
$ws.Cells.Item(316,1).Value = '16:52:47'
$ws.Cells.Item(316,2).Value = '16:53'
$ws.Cells.Item(316,3).Value = '10_OLMOS'
$ws.Cells.Item(316,4).Value = 1
$ws.Cells.Item(316,5).Value = 'LP1912'
$ws.Cells.Item(317,1).Value = '15:02:32'
$ws.Cells.Item(317,2).Value = '16:56'
$ws.Cells.Item(317,3).Value = '17_179 Y 38'
$ws.Cells.Item(317,4).Value = 114
$ws.Cells.Item(317,5).Value = 'LP1912'
$ws.Cells.Item(318,1).Value = '16:18:00'
$ws.Cells.Item(318,2).Value = '16:57'
$ws.Cells.Item(318,3).Value = '10_OLMOS'
$ws.Cells.Item(318,4).Value = 39
$ws.Cells.Item(318,5).Value = 'LP1912'
$ws.Cells.Item(319,1).Value = '16:52:47'
$ws.Cells.Item(319,2).Value = '17:04'
$ws.Cells.Item(319,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(319,4).Value = 12
$ws.Cells.Item(319,5).Value = 'LP1912'
$ws.Cells.Item(320,1).Value = '15:51:48'
$ws.Cells.Item(320,2).Value = '17:04'
$ws.Cells.Item(320,3).Value = '215A_EL PATO'
$ws.Cells.Item(320,4).Value = 73
$ws.Cells.Item(320,5).Value = 'LP1912'
$ws.Cells.Item(321,1).Value = '16:18:00'
$ws.Cells.Item(321,2).Value = '17:05'
$ws.Cells.Item(321,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(321,4).Value = 47
$ws.Cells.Item(321,5).Value = 'LP1912'
$ws.Cells.Item(322,1).Value = '16:18:00'
$ws.Cells.Item(322,2).Value = '17:05'
$ws.Cells.Item(322,3).Value = '215A_EL PATO'
$ws.Cells.Item(322,4).Value = 47
$ws.Cells.Item(322,5).Value = 'LP1912'
$ws.Cells.Item(323,1).Value = '16:18:00'
$ws.Cells.Item(323,2).Value = '17:05'
$ws.Cells.Item(323,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(323,4).Value = 47
$ws.Cells.Item(323,5).Value = 'LP1912'
$ws.Cells.Item(324,1).Value = '16:40:16'
$ws.Cells.Item(324,2).Value = '17:10'
$ws.Cells.Item(324,3).Value = '10_OLMOS'
$ws.Cells.Item(324,4).Value = 30
$ws.Cells.Item(324,5).Value = 'LP1912'
$ws.Cells.Item(325,1).Value = '16:52:47'
$ws.Cells.Item(325,2).Value = '17:16'
$ws.Cells.Item(325,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(325,4).Value = 24
$ws.Cells.Item(325,5).Value = 'LP1912'
$ws.Cells.Item(326,1).Value = '16:40:16'
$ws.Cells.Item(326,2).Value = '17:17'
$ws.Cells.Item(326,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(326,4).Value = 37
$ws.Cells.Item(326,5).Value = 'LP1912'
$ws.Cells.Item(327,1).Value = '16:52:47'
$ws.Cells.Item(327,2).Value = '17:20'
$ws.Cells.Item(327,3).Value = '16_SANTA ANA'
$ws.Cells.Item(327,4).Value = 28
$ws.Cells.Item(327,5).Value = 'LP1912'
$ws.Cells.Item(328,1).Value = '16:18:00'
$ws.Cells.Item(328,2).Value = '17:21'
$ws.Cells.Item(328,3).Value = '16_SANTA ANA'
$ws.Cells.Item(328,4).Value = 63
$ws.Cells.Item(328,5).Value = 'LP1912'
$ws.Cells.Item(329,1).Value = '15:51:48'
$ws.Cells.Item(329,2).Value = '17:21'
$ws.Cells.Item(329,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(329,4).Value = 90
$ws.Cells.Item(329,5).Value = 'LP1912'
$ws.Cells.Item(330,1).Value = '15:51:48'
$ws.Cells.Item(330,2).Value = '17:24'
$ws.Cells.Item(330,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(330,4).Value = 93
$ws.Cells.Item(330,5).Value = 'LP1912'
$ws.Cells.Item(331,1).Value = '15:51:48'
$ws.Cells.Item(331,2).Value = '17:28'
$ws.Cells.Item(331,3).Value = '14_ABASTO'
$ws.Cells.Item(331,4).Value = 97
$ws.Cells.Item(331,5).Value = 'LP1912'
$ws.Cells.Item(332,1).Value = '16:18:00'
$ws.Cells.Item(332,2).Value = '17:29'
$ws.Cells.Item(332,3).Value = '14_ABASTO'
$ws.Cells.Item(332,4).Value = 71
$ws.Cells.Item(332,5).Value = 'LP1912'
$ws.Cells.Item(333,1).Value = '16:52:47'
$ws.Cells.Item(333,2).Value = '17:31'
$ws.Cells.Item(333,3).Value = '15_ABASTO'
$ws.Cells.Item(333,4).Value = 39
$ws.Cells.Item(333,5).Value = 'LP1912'
$ws.Cells.Item(334,1).Value = '16:52:47'
$ws.Cells.Item(334,2).Value = '17:34'
$ws.Cells.Item(334,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(334,4).Value = 42
$ws.Cells.Item(334,5).Value = 'LP1912'
$ws.Cells.Item(335,1).Value = '16:40:16'
$ws.Cells.Item(335,2).Value = '17:35'
$ws.Cells.Item(335,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(335,4).Value = 55
$ws.Cells.Item(335,5).Value = 'LP1912'
$ws.Cells.Item(336,1).Value = '15:51:48'
$ws.Cells.Item(336,2).Value = '17:36'
$ws.Cells.Item(336,3).Value = '27_EL RETIRO'
$ws.Cells.Item(336,4).Value = 105
$ws.Cells.Item(336,5).Value = 'LP1912'
$ws.Cells.Item(337,1).Value = '16:40:16'
$ws.Cells.Item(337,2).Value = '17:37'
$ws.Cells.Item(337,3).Value = '27_EL RETIRO'
$ws.Cells.Item(337,4).Value = 57
$ws.Cells.Item(337,5).Value = 'LP1912'
$ws.Cells.Item(338,1).Value = '16:18:00'
$ws.Cells.Item(338,2).Value = '17:38'
$ws.Cells.Item(338,3).Value = '27_EL RETIRO'
$ws.Cells.Item(338,4).Value = 80
$ws.Cells.Item(338,5).Value = 'LP1912'
$ws.Cells.Item(339,1).Value = '15:51:48'
$ws.Cells.Item(339,2).Value = '17:38'
$ws.Cells.Item(339,3).Value = '17_ROMERO'
$ws.Cells.Item(339,4).Value = 107
$ws.Cells.Item(339,5).Value = 'LP1912'
$ws.Cells.Item(340,1).Value = '15:51:48'
$ws.Cells.Item(340,2).Value = '17:40'
$ws.Cells.Item(340,3).Value = '215B_EL PATO'
$ws.Cells.Item(340,4).Value = 109
$ws.Cells.Item(340,5).Value = 'LP1912'

# Append new rows
$ws.Cells.Item(341,1).Value = '16:40:16'
$ws.Cells.Item(341,2).Value = '17:41'
$ws.Cells.Item(341,3).Value = '16_SANTA ANA'
$ws.Cells.Item(341,4).Value = 61
$ws.Cells.Item(341,5).Value = 'LP1912'
$ws.Cells.Item(342,1).Value = '16:52:47'
$ws.Cells.Item(342,2).Value = '17:45'
$ws.Cells.Item(342,3).Value = '17_ROMERO'
$ws.Cells.Item(342,4).Value = 53
$ws.Cells.Item(342,5).Value = 'LP1912'
$ws.Cells.Item(343,1).Value = '15:51:48'
$ws.Cells.Item(343,2).Value = '17:50'
$ws.Cells.Item(343,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(343,4).Value = 119
$ws.Cells.Item(343,5).Value = 'LP1912'
$ws.Cells.Item(344,1).Value = '16:18:00'
$ws.Cells.Item(344,2).Value = '17:51'
$ws.Cells.Item(344,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(344,4).Value = 93
$ws.Cells.Item(344,5).Value = 'LP1912'
$ws.Cells.Item(345,1).Value = '16:18:00'
$ws.Cells.Item(345,2).Value = '17:52'
$ws.Cells.Item(345,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(345,4).Value = 94
$ws.Cells.Item(345,5).Value = 'LP1912'
$ws.Cells.Item(346,1).Value = '16:18:00'
$ws.Cells.Item(346,2).Value = '18:04'
$ws.Cells.Item(346,3).Value = '17_ROMERO'
$ws.Cells.Item(346,4).Value = 106
$ws.Cells.Item(346,5).Value = 'LP1912'
$ws.Cells.Item(347,1).Value = '16:40:16'
$ws.Cells.Item(347,2).Value = '18:21'
$ws.Cells.Item(347,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(347,4).Value = 101
$ws.Cells.Item(347,5).Value = 'LP1912'
$ws.Cells.Item(348,1).Value = '16:52:47'
$ws.Cells.Item(348,2).Value = '18:24'
$ws.Cells.Item(348,3).Value = '14_ABASTO'
$ws.Cells.Item(348,4).Value = 92
$ws.Cells.Item(348,5).Value = 'LP1912'
$ws.Cells.Item(349,1).Value = '16:52:47'
$ws.Cells.Item(349,2).Value = '18:27'
$ws.Cells.Item(349,3).Value = '215C_EL PATO'
$ws.Cells.Item(349,4).Value = 95
$ws.Cells.Item(349,5).Value = 'LP1912'
$ws.Cells.Item(350,1).Value = '16:40:16'
$ws.Cells.Item(350,2).Value = '18:28'
$ws.Cells.Item(350,3).Value = '215C_EL PATO'
$ws.Cells.Item(350,4).Value = 108
$ws.Cells.Item(350,5).Value = 'LP1912'
$ws.Cells.Item(351,1).Value = '16:40:16'
$ws.Cells.Item(351,2).Value = '18:32'
$ws.Cells.Item(351,3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(351,4).Value = 112
$ws.Cells.Item(351,5).Value = 'LP1912'
$ws.Cells.Item(352,1).Value = '16:52:47'
$ws.Cells.Item(352,2).Value = '18:48'
$ws.Cells.Item(352,3).Value = '14X44_ABASTO'
$ws.Cells.Item(352,4).Value = 116
$ws.Cells.Item(352,5).Value = 'LP1912'

# ===== Sheet 2: sheet2 =====
$ws = $wb.Worksheets.Item(2)

$ws.Range('A2').Value = 'Última actualización: 16:52:47'
$ws.Range('A3').Value = 'Total filas: 35'

# Update rows whose values changed due to re-sort
$ws.Cells.Item(39,1).Value = '16:52:47'
$ws.Cells.Item(39,2).Value = '18:27'
$ws.Cells.Item(39,3).Value = '215C_EL PATO'
$ws.Cells.Item(39,4).Value = 95
$ws.Cells.Item(39,5).Value = 'LP1912'

# Append new rows
$ws.Cells.Item(40,1).Value = '16:40:16'
$ws.Cells.Item(40,2).Value = '18:28'
$ws.Cells.Item(40,3).Value = '215C_EL PATO'
$ws.Cells.Item(40,4).Value = 108
$ws.Cells.Item(40,5).Value = 'LP1912'

# ===== Sheet 3: sheet3 =====
$ws = $wb.Worksheets.Item(3)

$ws.Range('A2').Value = 'Última actualización: 16:52:47'
$ws.Range('A3').Value = 'Total filas: 46'

# Update rows whose values changed due to re-sort
$ws.Cells.Item(50,1).Value = '16:52:47'
$ws.Cells.Item(50,2).Value = '18:03'
$ws.Cells.Item(50,3).Value = '215C_LA PLATA'
$ws.Cells.Item(50,4).Value = 71
$ws.Cells.Item(50,5).Value = 'L6203'

# Append new rows
$ws.Cells.Item(51,1).Value = '16:18:00'
$ws.Cells.Item(51,2).Value = '18:04'
$ws.Cells.Item(51,3).Value = '215C_LA PLATA'
$ws.Cells.Item(51,4).Value = 106
$ws.Cells.Item(51,5).Value = 'L6203'
